$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Texeira"
$ws.Range("C2").Value = "'02/05/2024"
$ws.Range("D2").Value = "31/05/2024"
$ws.Range("E2").Value = "internacional"
$ws.Range("F2").Value = "lazer"

# Row 3
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "patos"
$ws.Range("C3").Value = "'02/05/2024"
$ws.Range("D3").Value = "31/05/2024"
$ws.Range("E3").Value = "internacional"
$ws.Range("F3").Value = "lazer"
